$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "武器编号"
$ws.Range("B1").Value = "武器名称"
$ws.Range("C1").Value = "武器类型"
$ws.Range("D1").Value = "武器等级"
$ws.Range("E1").Value = "武器品质"
$ws.Range("F1").Value = "基础属性"
$ws.Range("G1").Value = "基础评分"
$ws.Range("H1").Value = "配方"

# Column F ("基础属性") is widened to a rendered width of 27 characters.
# The host subtracts a constant ~5/6-character padding when it stores the
# value typed into ColumnWidth, so we compensate here to land on an exact
# stored width of 27 (matches the saved <col .../> width).
$ws.Columns.Item(6).ColumnWidth = 26.16666666666667

# Mirror the author's final UI state: column H selected end-to-end.
$ws.Range("H1:H1048576").Select() | Out-Null
